# Update cryptos list figures (prices & 1h volume deltas) and swap
# the PancakeSwap / Decentraland row order, per the Apr 1 2023 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.492.74'
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').Value = '1.819.03'
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '''314.88'
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').Value = '''1.001'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').Value = '''0.5101'
$ws.Range('E7').Value = '  -6.72%  '
$ws.Range('D8').Value = '''0.3942'
$ws.Range('E8').Value = '  -2.36%  '
$ws.Range('D9').Value = '''0.08231'
$ws.Range('E9').Value = '  +7.00%  '
$ws.Range('D10').Value = '''1.108'
$ws.Range('E10').Value = '  -0.75%  '
$ws.Range('D11').Value = '''41.64'
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('D12').Value = '''21.11'
$ws.Range('E12').Value = '  +0.91%  '
$ws.Range('D13').Value = '''6.314'
$ws.Range('E13').Value = '  -0.19%  '
$ws.Range('D14').Value = '''1.001'
$ws.Range('E14').Value = '  +0.08%  '
$ws.Range('D15').Value = '''7.533'
$ws.Range('E15').Value = '  -1.70%  '
$ws.Range('D16').Value = '1.816.01'
$ws.Range('D17').Value = '''0.00001133'
$ws.Range('E17').Value = '  +4.67%  '
$ws.Range('D18').Value = '''92.45'
$ws.Range('E18').Value = '  +3.00%  '
$ws.Range('D19').Value = '''0.06651'
$ws.Range('E19').Value = '  +0.82%  '
$ws.Range('D20').Value = '''17.78'
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('D21').Value = '''0.9999'
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('D22').Value = '''6.087'
$ws.Range('E22').Value = '  +0.38%  '
$ws.Range('D23').Value = '28.531.18'
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('D24').Value = '''11.43'
$ws.Range('E24').Value = '  +2.84%  '
$ws.Range('D25').Value = '''2.267'
$ws.Range('E25').Value = '  +1.73%  '
$ws.Range('D26').Value = '''21.25'
$ws.Range('E26').Value = '  +2.39%  '
$ws.Range('D27').Value = '''155.77'
$ws.Range('E27').Value = '  -1.04%  '
$ws.Range('D28').Value = '2.026.25'
$ws.Range('E28').Value = '  -0.62%  '
$ws.Range('D29').Value = '''2.401'
$ws.Range('E29').Value = '  -2.92%  '
$ws.Range('D30').Value = '''125.60'
$ws.Range('E30').Value = '  +1.35%  '
$ws.Range('E31').Value = '  -1.46%  '
$ws.Range('D32').Value = '''0.1095'
$ws.Range('E32').Value = '  -1.30%  '
$ws.Range('D33').Value = '''5.776'
$ws.Range('E33').Value = '  +1.71%  '
$ws.Range('D34').Value = '''3.658'
$ws.Range('D35').Value = '''0.07079'
$ws.Range('E35').Value = '  -3.49%  '
$ws.Range('D36').Value = '''0.2220'
$ws.Range('E36').Value = '  -1.57%  '
$ws.Range('D37').Value = '''0.02344'
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('D38').Value = '''5.227'
$ws.Range('E38').Value = '  +0.41%  '
$ws.Range('D39').Value = '''8.820'
$ws.Range('E39').Value = '  -0.51%  '
$ws.Range('D40').Value = '''0.6306'
$ws.Range('E40').Value = '  +0.47%  '
$ws.Range('D41').Value = '''11.26'
$ws.Range('E41').Value = '  -0.97%  '
$ws.Range('D42').Value = '''1.179'
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('D43').Value = '''0.9996'
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('D44').Value = '''1.401'
$ws.Range('E44').Value = '  +0.53%  '
$ws.Range('D45').Value = '''13.53'
$ws.Range('E45').Value = '  -0.12%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').Value = '''3.732'
$ws.Range('E46').Value = '  +0.92%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '''0.5916'
$ws.Range('E47').Value = '  +1.11%  '
$ws.Range('D48').Value = '''124.86'
$ws.Range('E48').Value = '  -0.13%  '
$ws.Range('D49').Value = '''1.986'
$ws.Range('E49').Value = '  -0.94%  '
$ws.Range('D50').Value = '''1.184'
$ws.Range('E50').Value = '  -1.45%  '
$ws.Range('D51').Value = '''0.06889'
$ws.Range('E51').Value = '  +0.07%  '
